$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.961.39"
$ws.Range("E2").Value = "  -1.46%  "

# Row 3
$ws.Range("D3").Value = "2.413.39"
$ws.Range("E3").Value = "  -1.18%  "

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").Value = "'566.67"
$ws.Range("E5").Value = "  -1.93%  "

# Row 6
$ws.Range("D6").Value = "'138.93"
$ws.Range("E6").Value = "  -1.74%  "

# Row 7
$ws.Range("E7").Value = "  +0.28%  "

# Row 8
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = "  +1.82%  "

# Row 9
$ws.Range("D9").Value = "2.396.32"
$ws.Range("E9").Value = "  -1.70%  "

# Row 10
$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -2.57%  "

# Row 11
$ws.Range("E11").Value = "  -0.52%  "

# Row 12
$ws.Range("D12").Value = "'5.05"
$ws.Range("E12").Value = "  -2.28%  "

# Row 13
$ws.Range("D13").Value = "'0.335"
$ws.Range("E13").Value = "  -1.16%  "

# Row 14
$ws.Range("D14").Value = "'25.85"
$ws.Range("E14").Value = "  -0.15%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.847.18"
$ws.Range("E15").Value = "  -1.60%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000169"
$ws.Range("E16").Value = "  -2.07%  "

# Row 17
$ws.Range("D17").Value = "61.002.96"
$ws.Range("E17").Value = "  -1.34%  "

# Row 18
$ws.Range("D18").Value = "2.397.35"
$ws.Range("E18").Value = "  -1.95%  "

# Row 19
$ws.Range("D19").Value = "'8.15"
$ws.Range("E19").Value = "  +13.15%  "

# Row 20
$ws.Range("D20").Value = "'10.49"
$ws.Range("E20").Value = "  -1.23%  "

# Row 21
$ws.Range("D21").Value = "'322.53"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22
$ws.Range("D22").Value = "'4.04"
$ws.Range("E22").Value = "  -0.33%  "

# Row 23
$ws.Range("D23").Value = "'6.16"
$ws.Range("E23").Value = "  +3.51%  "

# Row 25
$ws.Range("D25").Value = "'1.82"
$ws.Range("E25").Value = "  -4.27%  "

# Row 26
$ws.Range("D26").Value = "'64.28"
$ws.Range("E26").Value = "  -1.17%  "

# Row 27
$ws.Range("D27").Value = "'573.37"
$ws.Range("E27").Value = "  -1.93%  "

# Row 28
$ws.Range("D28").Value = "'8.19"
$ws.Range("E28").Value = "  -10.50%  "

# Row 29
$ws.Range("D29").Value = "2.531.85"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0921"
$ws.Range("E30").Value = "  -1.85%  "

# Row 31
$ws.Range("D31").Value = "'7.93"
$ws.Range("E31").Value = "  +0.78%  "

# Row 32
$ws.Range("D32").Value = "'1.32"
$ws.Range("E32").Value = "  -4.16%  "

# Row 33
$ws.Range("D33").Value = "'1.81"
$ws.Range("E33").Value = "  -3.47%  "

# Row 34
$ws.Range("E34").Value = "  -0.94%  "

# Row 35
$ws.Range("E35").Value = "  +0.19%  "

# Row 36
$ws.Range("E36").Value = "  +0.56%  "

# Row 37
$ws.Range("D37").Value = "'152.41"
$ws.Range("E37").Value = "  -0.25%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.368"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'4.54"
$ws.Range("E39").Value = "  -4.74%  "

# Row 40
$ws.Range("D40").Value = "'18.17"
$ws.Range("E40").Value = "  -0.78%  "

# Row 41
$ws.Range("D41").Value = "'5.10"
$ws.Range("E41").Value = "  -1.18%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").Value = "'1.65"
$ws.Range("E43").Value = "  -1.32%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.37"
$ws.Range("E44").Value = "  -0.45%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0294"
$ws.Range("E45").Value = "  +7.83%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'142.50"
$ws.Range("E46").Value = "  +1.32%  "

# Row 47
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'3.51"
$ws.Range("E47").Value = "  -1.64%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.584"
$ws.Range("E48").Value = "  -2.39%  "

# Row 49
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0499"
$ws.Range("E49").Value = "  -2.27%  "

# Row 50
$ws.Range("D50").Value = "'19.15"
$ws.Range("E50").Value = "  -2.28%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0901"
$ws.Range("E51").Value = "  +0.56%  "
